$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new "self-grade" column (C) alongside the existing "Points" column (B).
$ws.Range("C3").Value  = 1.5
$ws.Range("C4").Value  = 2
$ws.Range("C5").Value  = 1
$ws.Range("C6").Value  = 2
$ws.Range("C7").Value  = 1
$ws.Range("C8").Value  = 1
$ws.Range("C9").Value  = 1
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 2
$ws.Range("C13").Value = 1

# Total row: sum the new column, mirroring the existing B14 total formula.
$ws.Range("C14").Formula = "=SUM(C3:C13)"

# Leave the selection where the editor finished working, in the newly-added column.
$ws.Range("D14").Select()
